$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Selplg"
$ws.Range("C2").Value = "Sell"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.156459
$ws.Range("H2").Value = 6.469377
$ws.Range("I2").Value = 0.01389412936885011
$ws.Range("J2").Value = 0.01392987523772938
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.146719
$ws.Range("N2").Value = 0.440157
$ws.Range("O2").Value = 0.002213357657235064
$ws.Range("P2").Value = 0.002213357657235064
$ws.Range("Q2").Value = 0.316393508021
$ws.Range("R2").Value = 2.847541572189
$ws.Range("S2").Value = 0.00003075267762915897
$ws.Range("T2").Value = 0.00003083179602175742

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Selplg"
$ws.Range("C3").Value = "Sell"
$ws.Range("D3").Value = "M1"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.156459
$ws.Range("H3").Value = 6.469377
$ws.Range("I3").Value = 0.01389412936885011
$ws.Range("J3").Value = 0.01392987523772938
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 28.93198366666667
$ws.Range("N3").Value = 86.795951
$ws.Range("O3").Value = 0.4364589970461662
$ws.Range("P3").Value = 0.4364589970461662
$ws.Range("Q3").Value = 62.39063656583633
$ws.Range("R3").Value = 561.515729092527
$ws.Range("S3").Value = 0.006064217769158001
$ws.Range("T3").Value = 0.006079819375237589

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Selplg"
$ws.Range("C4").Value = "Sell"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.156459
$ws.Range("H4").Value = 6.469377
$ws.Range("I4").Value = 0.01389412936885011
$ws.Range("J4").Value = 0.01392987523772938
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 37.20927366666667
$ws.Range("N4").Value = 111.627821
$ws.Range("O4").Value = 0.5613276452965987
$ws.Range("P4").Value = 0.5613276452965988
$ws.Range("Q4").Value = 80.24027308194633
$ws.Range("R4").Value = 722.162457737517
$ws.Range("S4").Value = 0.007799158922062949
$ws.Range("T4").Value = 0.00781922406647003

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Selplg"
$ws.Range("C5").Value = "Sell"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.009804666666666
$ws.Range("H5").Value = 9.029413999999999
$ws.Range("I5").Value = 0.01939226083762105
$ws.Range("J5").Value = 0.01944215192433629
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.146719
$ws.Range("N5").Value = 0.440157
$ws.Range("O5").Value = 0.002213357657235064
$ws.Range("P5").Value = 0.002213357657235064
$ws.Range("Q5").Value = 0.4415955308886667
$ws.Range("R5").Value = 3.974359777998
$ws.Range("S5").Value = 0.0000429220090160482
$ws.Range("T5").Value = 0.00004303243583485716

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Selplg"
$ws.Range("C6").Value = "Sell"
$ws.Range("D6").Value = "M1"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.009804666666666
$ws.Range("H6").Value = 9.029413999999999
$ws.Range("I6").Value = 0.01939226083762105
$ws.Range("J6").Value = 0.01944215192433629
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 28.93198366666667
$ws.Range("N6").Value = 86.795951
$ws.Range("O6").Value = 0.4364589970461662
$ws.Range("P6").Value = 0.4364589970461662
$ws.Range("Q6").Value = 87.0796194558571
$ws.Range("R6").Value = 783.7165751027139
$ws.Range("S6").Value = 0.008463926715645729
$ws.Range("T6").Value = 0.008485702129315008

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Selplg"
$ws.Range("C7").Value = "Sell"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.009804666666666
$ws.Range("H7").Value = 9.029413999999999
$ws.Range("I7").Value = 0.01939226083762105
$ws.Range("J7").Value = 0.01944215192433629
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 37.20927366666667
$ws.Range("N7").Value = 111.627821
$ws.Range("O7").Value = 0.5613276452965987
$ws.Range("P7").Value = 0.5613276452965988
$ws.Range("Q7").Value = 111.9926455252104
$ws.Range("R7").Value = 1007.933809726894
$ws.Range("S7").Value = 0.01088541211295927
$ws.Range("T7").Value = 0.01091341735918643

# Row 8
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Selplg"
$ws.Range("C8").Value = "Sell"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 80.649016
$ws.Range("H8").Value = 241.947048
$ws.Range("I8").Value = 0.519624004803459
$ws.Range("J8").Value = 0.5209608580203196
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.146719
$ws.Range("N8").Value = 0.440157
$ws.Range("O8").Value = 0.002213357657235064
$ws.Range("P8").Value = 0.002213357657235064
$ws.Range("Q8").Value = 11.832742978504
$ws.Range("R8").Value = 106.494686806536
$ws.Range("S8").Value = 0.001150113769914886
$ws.Range("T8").Value = 0.001153072704219023

# Row 9
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Selplg"
$ws.Range("C9").Value = "Sell"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 80.649016
$ws.Range("H9").Value = 241.947048
$ws.Range("I9").Value = 0.519624004803459
$ws.Range("J9").Value = 0.5209608580203196
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 28.93198366666667
$ws.Range("N9").Value = 86.795951
$ws.Range("O9").Value = 0.4364589970461662
$ws.Range("P9").Value = 0.4364589970461662
$ws.Range("Q9").Value = 2333.336013644739
$ws.Range("R9").Value = 21000.02412280265
$ws.Range("S9").Value = 0.22679457197763
$ws.Range("T9").Value = 0.2273780535918589

# Row 10
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Selplg"
$ws.Range("C10").Value = "Sell"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 80.649016
$ws.Range("H10").Value = 241.947048
$ws.Range("I10").Value = 0.519624004803459
$ws.Range("J10").Value = 0.5209608580203196
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 37.20927366666667
$ws.Range("N10").Value = 111.627821
$ws.Range("O10").Value = 0.5613276452965987
$ws.Range("P10").Value = 0.5613276452965988
$ws.Range("Q10").Value = 3000.891307291379
$ws.Range("R10").Value = 27008.02176562241
$ws.Range("S10").Value = 0.2916793190559142
$ws.Range("T10").Value = 0.2924297317242417

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Selplg"
$ws.Range("C11").Value = "Sell"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 68.19636666666668
$ws.Range("H11").Value = 204.5891
$ws.Range("I11").Value = 0.4393912153916231
$ws.Range("J11").Value = 0.4405216511573432
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.146719
$ws.Range("N11").Value = 0.440157
$ws.Range("O11").Value = 0.002213357657235064
$ws.Range("P11").Value = 0.002213357657235064
$ws.Range("Q11").Value = 10.00570272096667
$ws.Range("R11").Value = 90.05132448870002
$ws.Range("S11").Value = 0.0009725299111088701
$ws.Range("T11").Value = 0.0009750319697669392

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Selplg"
$ws.Range("C12").Value = "Sell"
$ws.Range("D12").Value = "M1"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 68.19636666666668
$ws.Range("H12").Value = 204.5891
$ws.Range("I12").Value = 0.4393912153916231
$ws.Range("J12").Value = 0.4405216511573432
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 28.93198366666667
$ws.Range("N12").Value = 86.795951
$ws.Range("O12").Value = 0.4364589970461662
$ws.Range("P12").Value = 0.4364589970461662
$ws.Range("Q12").Value = 1973.056166526012
$ws.Range("R12").Value = 17757.5054987341
$ws.Range("S12").Value = 0.1917762491807238
$ws.Range("T12").Value = 0.1922696380412551

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Selplg"
$ws.Range("C13").Value = "Sell"
$ws.Range("D13").Value = "M2"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 68.19636666666668
$ws.Range("H13").Value = 204.5891
$ws.Range("I13").Value = 0.4393912153916231
$ws.Range("J13").Value = 0.4405216511573432
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 37.20927366666667
$ws.Range("N13").Value = 111.627821
$ws.Range("O13").Value = 0.5613276452965987
$ws.Range("P13").Value = 0.5613276452965988
$ws.Range("Q13").Value = 2537.537270372345
$ws.Range("R13").Value = 22837.8354333511
$ws.Range("S13").Value = 0.2466424362997904
$ws.Range("T13").Value = 0.2472769811463212

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Selplg"
$ws.Range("C14").Value = "Sell"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.19484
$ws.Range("H14").Value = 2.38968
$ws.Range("I14").Value = 0.007698389598446743
$ws.Range("J14").Value = 0.005145463660271636
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.146719
$ws.Range("N14").Value = 0.440157
$ws.Range("O14").Value = 0.002213357657235064
$ws.Range("P14").Value = 0.002213357657235064
$ws.Range("Q14").Value = 0.17530572996
$ws.Range("R14").Value = 1.05183437976
$ws.Range("S14").Value = 0.00001703928956610087
$ws.Range("T14").Value = 0.00001138875139248698

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Selplg"
$ws.Range("C15").Value = "Sell"
$ws.Range("D15").Value = "M1"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.19484
$ws.Range("H15").Value = 2.38968
$ws.Range("I15").Value = 0.007698389598446743
$ws.Range("J15").Value = 0.005145463660271636
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 28.93198366666667
$ws.Range("N15").Value = 86.795951
$ws.Range("O15").Value = 0.4364589970461662
$ws.Range("P15").Value = 0.4364589970461662
$ws.Range("Q15").Value = 34.56909136428001
$ws.Range("R15").Value = 207.41454818568
$ws.Range("S15").Value = 0.003360031403008704
$ws.Range("T15").Value = 0.002245783908499653

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Selplg"
$ws.Range("C16").Value = "Sell"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.19484
$ws.Range("H16").Value = 2.38968
$ws.Range("I16").Value = 0.007698389598446743
$ws.Range("J16").Value = 0.005145463660271636
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 37.20927366666667
$ws.Range("N16").Value = 111.627821
$ws.Range("O16").Value = 0.5613276452965987
$ws.Range("P16").Value = 0.5613276452965988
$ws.Range("Q16").Value = 44.45912854788001
$ws.Range("R16").Value = 266.75477128728
$ws.Range("S16").Value = 0.004321318905871938
$ws.Range("T16").Value = 0.002888291000379496

Write-Host "Done"